$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'31.182.04"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "'1.990.61"
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'0.7933"
$ws.Range("E5").Value = "  +67.33%  "
$ws.Range("E6").Value = "  +3.41%  "
$ws.Range("D7").Value = "'0.9994"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.3495"
$ws.Range("E8").Value = "  +20.75%  "
$ws.Range("D9").Value = "'28.13"
$ws.Range("E9").Value = "  +25.85%  "
$ws.Range("E10").Value = "  +6.81%  "
$ws.Range("D11").Value = "'0.8452"
$ws.Range("E11").Value = "  +8.88%  "
$ws.Range("D12").Value = "'0.08178"
$ws.Range("E12").Value = "  +4.64%  "
$ws.Range("D13").Value = "'100.38"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").Value = "'1.988.61"
$ws.Range("E14").Value = "  +5.55%  "
$ws.Range("D15").Value = "'5.627"
$ws.Range("E15").Value = "  +6.98%  "
$ws.Range("D16").Value = "'15.42"
$ws.Range("E16").Value = "  +16.64%  "
$ws.Range("D17").Value = "'272.94"
$ws.Range("E17").Value = "  -4.18%  "
$ws.Range("D18").Value = "'31.182.30"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("D19").Value = "'5.873"
$ws.Range("E19").Value = "  +9.48%  "
$ws.Range("D20").Value = "'0.000007938"
$ws.Range("E20").Value = "  +5.33%  "
$ws.Range("D21").Value = "'2.248.93"
$ws.Range("E21").Value = "  +5.71%  "
$ws.Range("D22").Value = "'0.9989"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'7.049"
$ws.Range("E24").Value = "  +8.99%  "
$ws.Range("D25").Value = "'9.984"
$ws.Range("D26").Value = "'0.1507"
$ws.Range("E26").Value = "  +55.22%  "
$ws.Range("D27").Value = "'165.57"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").Value = "'19.88"
$ws.Range("E28").Value = "  +3.87%  "
$ws.Range("D29").Value = "'2.343"
$ws.Range("E29").Value = "  +22.29%  "
$ws.Range("D30").Value = "'1.595"
$ws.Range("E30").Value = "  +6.21%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'1.355"
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.583"
$ws.Range("E32").Value = "  +7.79%  "
$ws.Range("D33").Value = "'4.413"
$ws.Range("E33").Value = "  +5.24%  "
$ws.Range("D34").Value = "'0.05258"
$ws.Range("E34").Value = "  +8.38%  "
$ws.Range("D35").Value = "'0.7796"
$ws.Range("E35").Value = "  +11.47%  "
$ws.Range("D36").Value = "'1.215"
$ws.Range("E36").Value = "  +7.42%  "
$ws.Range("D37").Value = "'2.760"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "'0.9978"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").Value = "'0.02007"
$ws.Range("E39").Value = "  +4.62%  "
$ws.Range("D40").Value = "'2.899"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "'6.646"
$ws.Range("E41").Value = "  +5.52%  "
$ws.Range("D42").Value = "'79.57"
$ws.Range("E42").Value = "  +4.58%  "
$ws.Range("D43").Value = "'0.4654"
$ws.Range("E43").Value = "  +9.30%  "
$ws.Range("D44").Value = "'2.124"
$ws.Range("E44").Value = "  +6.63%  "
$ws.Range("D45").Value = "'0.8530"
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("D46").Value = "'104.62"
$ws.Range("E46").Value = "  +2.95%  "
$ws.Range("D47").Value = "'0.9989"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "'7.671"
$ws.Range("E48").Value = "  +9.11%  "
$ws.Range("D49").Value = "'9.875"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("D50").Value = "'36.86"
$ws.Range("E50").Value = "  +4.63%  "
$ws.Range("E51").Value = "  +8.20%  "
